$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the title text (column A) between row 3 and row 4 ---
# NOTE: must read via Value() (with parens) -- the bare `.Value` getter
# (no parens) does not invoke the accessor in this host and instead yields
# a stringified member-signature, so always call it like a method when reading.
$titleRow3 = $ws.Range("A3").Value()
$titleRow4 = $ws.Range("A4").Value()
$ws.Range("A3").Value = $titleRow4
$ws.Range("A4").Value = $titleRow3

# --- Swap the displayed uri text (column E) between row 3 and row 4 ---
$uriRow3 = $ws.Range("E3").Value()
$uriRow4 = $ws.Range("E4").Value()
$ws.Range("E3").Value = $uriRow4
$ws.Range("E4").Value = $uriRow3

# --- Move the "#1960" hyperlink sub-address (location) from the E3 link to the E4 link ---
# (materialize the collection into an array first -- indexing via .Item() directly on the
#  live collection does not resolve property reads/writes correctly in this host)
$links = @($ws.Hyperlinks)
$linkE3 = $links[1]
$linkE4 = $links[2]

$linkE3.SubAddress = ""
$linkE4.SubAddress = "1960"

Write-Host "done"
